$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("I9").Value = 283
$ws.Range("J9").Value = 111
$ws.Range("K9").Value = 283
$ws.Range("L9").Value = 111
$ws.Range("M9").Value = -114
$ws.Range("N9").Value = -449

$ws.Range("H15").Value = 1213.7174
$ws.Range("I15").Value = 1213.7174
$ws.Range("K15").Value = 3641.1522
$ws.Range("M15").Value = -3472.1522

$ws.Range("H17").Value = 1014.0238
$ws.Range("J17").Value = 1069.5897
$ws.Range("L17").Value = 3208.7691
$ws.Range("N17").Value = -3544.7691

$ws.Range("H28").Value = 32582.334
$ws.Range("I28").Value = 40758.56
$ws.Range("K28").Value = 40758.56
$ws.Range("M28").Value = -40273.56

$ws.Range("H51").Value = 5365.933
$ws.Range("I51").Value = 3500
$ws.Range("J51").Value = 7498.4287
$ws.Range("K51").Value = 3500
$ws.Range("L51").Value = 7498.4287
$ws.Range("M51").Value = -3016
$ws.Range("N51").Value = -8466.4287

$ws.Range("H88").Value = 1750
$ws.Range("J88").Value = 1500
$ws.Range("L88").Value = 1500
$ws.Range("N88").Value = -2312

$ws.Range("H91").Value = 1750
$ws.Range("J91").Value = 1500
$ws.Range("L91").Value = 1500
$ws.Range("N91").Value = -4308

$ws.Range("H138").Value = 2345.8286
$ws.Range("I138").Value = 1861.0714
$ws.Range("J138").Value = 4284.857
$ws.Range("K138").Value = 5583.2142
$ws.Range("L138").Value = 12854.571
$ws.Range("M138").Value = -443.2142000000003
$ws.Range("N138").Value = -23134.571

$ws.Range("H141").Value = 1922.65
$ws.Range("I141").Value = 1134
$ws.Range("J141").Value = 2886.5557
$ws.Range("K141").Value = 3402
$ws.Range("L141").Value = 8659.667099999999
$ws.Range("M141").Value = 1778
$ws.Range("N141").Value = -19019.6671

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 146620.17
$ws.Range("I32").Value = 146620.17
$ws.Range("K32").Value = 146620.17
$ws.Range("M32").Value = -146333.17

$ws.Range("H61").Value = 2351.25
$ws.Range("I61").Value = 2351.25
$ws.Range("K61").Value = 2351.25
$ws.Range("M61").Value = -2139.25

$ws.Range("H74").Value = 3365.9092
$ws.Range("I74").Value = 3683.2
$ws.Range("K74").Value = 3683.2
$ws.Range("M74").Value = -2809.2

$ws.Range("H77").Value = 3365.9092
$ws.Range("I77").Value = 3683.2
$ws.Range("K77").Value = 18416
$ws.Range("M77").Value = -14048

$ws.Range("H122").Value = 20836934
$ws.Range("I122").Value = 41669132
$ws.Range("K122").Value = 125007396
$ws.Range("M122").Value = -125004946

$ws.Range("H136").Value = 2351.25
$ws.Range("I136").Value = 2351.25
$ws.Range("K136").Value = 7053.75
$ws.Range("M136").Value = -4503.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 0
$ws.Range("N6").ClearContents()

$ws.Range("H107").Value = 10007259
$ws.Range("I107").Value = 3689.8647
$ws.Range("K107").Value = 3689.8647
$ws.Range("M107").Value = -1769.8647

$ws.Range("H134").Value = 2730.926
$ws.Range("I134").Value = 2709.4
$ws.Range("J134").Value = 3000
$ws.Range("K134").Value = 8128.200000000001
$ws.Range("L134").Value = 9000
$ws.Range("M134").Value = -5593.200000000001
$ws.Range("N134").Value = -14070

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 52.555557
$ws.Range("I7").Value = 46.142857
$ws.Range("K7").Value = 46.142857
$ws.Range("M7").Value = 66.85714300000001

$ws.Range("H31").Value = 3440.2856
$ws.Range("I31").Value = 2898.6
$ws.Range("J31").Value = 3609.5625
$ws.Range("K31").Value = 2898.6
$ws.Range("L31").Value = 3609.5625
$ws.Range("M31").Value = -2603.6
$ws.Range("N31").Value = -4199.5625

$ws.Range("H34").Value = 3440.2856
$ws.Range("I34").Value = 2898.6
$ws.Range("J34").Value = 3609.5625
$ws.Range("K34").Value = 2898.6
$ws.Range("L34").Value = 3609.5625
$ws.Range("M34").Value = -2696.6
$ws.Range("N34").Value = -4013.5625

$ws.Range("H48").Value = 0
$ws.Range("J48").Value = 0
$ws.Range("N48").ClearContents()

$ws.Range("H86").Value = 10000
$ws.Range("I86").Value = 10000
$ws.Range("K86").Value = 10000
$ws.Range("M86").Value = -8877

$ws.Range("H89").Value = 10000
$ws.Range("I89").Value = 10000
$ws.Range("K89").Value = 50000
$ws.Range("M89").Value = -44384

$ws.Range("H105").Value = 1385
$ws.Range("I105").Value = 1390.1333
$ws.Range("K105").Value = 1390.1333
$ws.Range("M105").Value = 356.8667

$ws.Range("H134").Value = 3346.3333
$ws.Range("I134").Value = 1813.5714
$ws.Range("J134").Value = 4687.5
$ws.Range("K134").Value = 5440.7142
$ws.Range("L134").Value = 14062.5
$ws.Range("M134").Value = -2905.7142
$ws.Range("N134").Value = -19132.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 798.82355
$ws.Range("I5").Value = 341.45456
$ws.Range("J5").Value = 1637.3334
$ws.Range("K5").Value = 1024.36368
$ws.Range("L5").Value = 4912.0002
$ws.Range("M5").Value = -912.3636799999999
$ws.Range("N5").Value = -5136.0002

$ws.Range("H9").Value = 14405.857
$ws.Range("I9").Value = 14405.857
$ws.Range("K9").Value = 43217.571
$ws.Range("M9").Value = -42993.571

$ws.Range("H62").Value = 3961.1667
$ws.Range("J62").Value = 3255.6667
$ws.Range("L62").Value = 9767.000100000001
$ws.Range("N62").Value = -11139.0001

$ws.Range("H65").Value = 3961.1667
$ws.Range("J65").Value = 3255.6667
$ws.Range("L65").Value = 29301.0003
$ws.Range("N65").Value = -36165.0003

$ws.Range("H68").Value = 2022.7693
$ws.Range("I68").Value = 1549.5
$ws.Range("K68").Value = 4648.5
$ws.Range("M68").Value = -3837.5

$ws.Range("H71").Value = 2022.7693
$ws.Range("I71").Value = 1549.5
$ws.Range("K71").Value = 13945.5
$ws.Range("M71").Value = -9889.5

$ws.Range("H86").Value = 670.6923
$ws.Range("I86").Value = 620.8333
$ws.Range("K86").Value = 1862.4999
$ws.Range("M86").Value = -676.4999

$ws.Range("H89").Value = 670.6923
$ws.Range("I89").Value = 620.8333
$ws.Range("K89").Value = 5587.4997
$ws.Range("M89").Value = 340.5002999999997

$ws.Range("H107").Value = 1412.2174
$ws.Range("I107").Value = 1431.5
$ws.Range("K107").Value = 4294.5
$ws.Range("M107").Value = -2374.5

$ws.Range("H113").Value = 1322.2307
$ws.Range("J113").Value = 1599
$ws.Range("L113").Value = 4797
$ws.Range("N113").Value = -9137

$ws.Range("H135").Value = 798.82355
$ws.Range("I135").Value = 341.45456
$ws.Range("J135").Value = 1637.3334
$ws.Range("K135").Value = 3073.09104
$ws.Range("L135").Value = 14736.0006
$ws.Range("M135").Value = -538.0910400000002
$ws.Range("N135").Value = -19806.0006

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 41669588
$ws.Range("I122").Value = 2784
$ws.Range("K122").Value = 8352
$ws.Range("M122").Value = -5902

$ws.Range("H138").Value = 59999
$ws.Range("J138").Value = 59999
$ws.Range("L138").Value = 59999
$ws.Range("N138").Value = -70279

$ws.Range("H139").Value = 86000
$ws.Range("J139").Value = 86000
$ws.Range("L139").Value = 86000
$ws.Range("N139").Value = -96280

$ws.Range("H141").Value = 65071.332
$ws.Range("I141").Value = 20000
$ws.Range("J141").Value = 87607
$ws.Range("K141").Value = 20000
$ws.Range("L141").Value = 87607
$ws.Range("M141").Value = -14820
$ws.Range("N141").Value = -97967

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3254.6538
$ws.Range("I40").Value = 2906.3
$ws.Range("J40").Value = 4415.8335
$ws.Range("K40").Value = 2906.3
$ws.Range("L40").Value = 4415.8335
$ws.Range("M40").Value = -2770.3
$ws.Range("N40").Value = -4687.8335

$ws.Range("H46").Value = 2041.1428
$ws.Range("J46").Value = 2055.2693
$ws.Range("L46").Value = 2055.2693
$ws.Range("N46").Value = -2431.2693

$ws.Range("H132").Value = 6632.6665
$ws.Range("I132").Value = 3299.6667
$ws.Range("J132").Value = 9132.416999999999
$ws.Range("K132").Value = 9899.000100000001
$ws.Range("L132").Value = 27397.251
$ws.Range("M132").Value = -7369.000100000001
$ws.Range("N132").Value = -32457.251

$ws.Range("H136").Value = 6591.8
$ws.Range("I136").Value = 4252.737
$ws.Range("K136").Value = 12758.211
$ws.Range("M136").Value = -10208.211

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 771.6
$ws.Range("I107").Value = 853
$ws.Range("J107").Value = 649.5
$ws.Range("K107").Value = 2559
$ws.Range("L107").Value = 1948.5
$ws.Range("M107").Value = -639
$ws.Range("N107").Value = -5788.5

$ws.Range("H113").Value = 1309.7273
$ws.Range("I113").Value = 1637.3334
$ws.Range("K113").Value = 4912.0002
$ws.Range("M113").Value = -2742.0002

$ws.Range("H132").Value = 875420.5
$ws.Range("I132").Value = 1544574.5
$ws.Range("J132").Value = 5520.4
$ws.Range("K132").Value = 4633723.5
$ws.Range("L132").Value = 16561.2
$ws.Range("M132").Value = -4631193.5
$ws.Range("N132").Value = -21621.2

$ws.Range("H136").Value = 7833.3335
$ws.Range("I136").Value = 3833.3333
$ws.Range("J136").Value = 11833.333
$ws.Range("K136").Value = 11499.9999
$ws.Range("L136").Value = 35499.999
$ws.Range("M136").Value = -8949.999899999999
$ws.Range("N136").Value = -40599.999

$ws.Range("H140").Value = 55000
$ws.Range("I140").Value = 55000
$ws.Range("J140").Value = 40771
$ws.Range("K140").Value = 55000
$ws.Range("M140").Value = -49820
$ws.Range("N140").ClearContents()
